$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.899.12"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.635.27"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'216.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").Value = "'0.5076"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.2579"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("D9").Value = "'0.06363"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "'19.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("D11").Value = "'0.07754"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.674.85"
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.262"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").Value = "'0.5516"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").Value = "0.0₅7700"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "25.920.33"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "'1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "'194.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").Value = "'9.915"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "'6.050"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "'1.906"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").Value = "'142.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").Value = "'0.1239"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.17%  "
$ws.Range("D27").Value = "'6.816"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.87%  "
$ws.Range("D28").Value = "'15.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("D29").Value = "'1.243"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").Value = "'0.04885"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").Value = "'3.254"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "'3.193"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").Value = "'1.543"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").Value = "'2.380"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").Value = "'0.9061"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.33%  "
$ws.Range("D36").Value = "'2.560"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("D37").Value = "'0.5499"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("D38").Value = "1.123.48"
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("D39").Value = "'0.01557"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").Value = "'5.572"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "'0.8043"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.83%  "
$ws.Range("D43").Value = "'97.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("E44").Value = "  -4.14%  "
$ws.Range("D45").Value = "1.773.42"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").Value = "'0.4451"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.85%  "
$ws.Range("D47").Value = "'54.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("D50").Value = "'7.551"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.61%  "
$ws.Range("E51").Value = "  -0.28%  "
